# Update "想去人数" (want-to-go count) values in both the "展览" sheet
# and the "全部类型" sheet, which contain duplicated data for the same
# events.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 460
    $ws.Range("F4").Value = 20
}
